$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.011.91'
$ws.Range("E2").Value = '  +7.12%  '

$ws.Range("D3").Value = '3.022.07'
$ws.Range("E3").Value = '  +4.31%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.00%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.42'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").Value = '3.018.92'
$ws.Range("E8").Value = '  +4.25%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.519'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.75%  '

$ws.Range("E10").Value = '  +0.11%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.155'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.454'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.35%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000248'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.01%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.48'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +8.11%  '

$ws.Range("E15").Value = '  +0.76%  '

$ws.Range("D16").Value = '66.004.56'
$ws.Range("E16").Value = '  +7.13%  '

$ws.Range("D17").Value = '3.519.02'
$ws.Range("E17").Value = '  +4.20%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.98'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.40%  '

$ws.Range("D19").Value = '3.018.34'
$ws.Range("E19").Value = '  +4.25%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '462.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.65%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.83'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.67%  '

$ws.Range("E22").Value = '  +4.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.39'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.25%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.17'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.42%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.60'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.15%  '

$ws.Range("E26").Value = '  +12.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +9.67%  '

$ws.Range("E28").Value = '  -0.06%  '

$ws.Range("B29").Value = 'NEARProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.98'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +14.00%  '

$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.42'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +18.62%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0000105'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.70%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.62'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.46%  '

$ws.Range("E33").Value = '  +5.24%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.06'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.80%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.995'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.56%  '

$ws.Range("B37").Value = 'Stacks'
$ws.Range("C37").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.20'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +13.60%  '

$ws.Range("B38").Value = 'Filecoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.80'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.48%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.05'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.24%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.44'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.23%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '45.10'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +14.13%  '

$ws.Range("E42").Value = '  +7.57%  '

$ws.Range("E43").Value = '  +13.53%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.51'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.56%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '392.88'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +13.05%  '

$ws.Range("D46").Value = '2.801.77'
$ws.Range("E46").Value = '  +4.13%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0355'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.94%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '134.95'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.36%  '

$ws.Range("E49").Value = '  -0.09%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.80'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +10.07%  '

$ws.Range("E51").Value = '  +3.96%  '
